# Each (Região) block of 21 rows (Brasil: 2-22, Nordeste: 23-43, Sergipe: 44-64)
# is shifted up by one quarter: the "Trimestre" (C) and "Valor" (D) columns
# both advance by one row, a new trailing quarter/value is appended, and the
# row that no longer has a following value loses its D cell entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; C = "01/04/2019"; D = 12.1 },
    @{ Row = 3; C = "01/07/2019"; D = 11.9 },
    @{ Row = 4; C = "01/10/2019"; D = 11.1 },
    @{ Row = 5; C = "01/01/2020"; D = 12.4 },
    @{ Row = 6; C = "01/04/2020"; D = 13.6 },
    @{ Row = 7; C = "01/07/2020"; D = 14.9 },
    @{ Row = 8; C = "01/10/2020"; D = 14.2 },
    @{ Row = 9; C = "01/01/2021"; D = 14.9 },
    @{ Row = 10; C = "01/04/2021"; D = 14.2 },
    @{ Row = 11; C = "01/07/2021"; D = 12.6 },
    @{ Row = 12; C = "01/10/2021"; D = 11.1 },
    @{ Row = 13; C = "01/01/2022"; D = 11.1 },
    @{ Row = 14; C = "01/04/2022"; D = 9.300000000000001 },
    @{ Row = 15; C = "01/07/2022"; D = 8.699999999999999 },
    @{ Row = 16; C = "01/10/2022"; D = 7.9 },
    @{ Row = 17; C = "01/01/2023"; D = 8.800000000000001 },
    @{ Row = 18; C = "01/04/2023"; D = 8 },
    @{ Row = 19; C = "01/07/2023"; D = 7.7 },
    @{ Row = 20; C = "01/10/2023"; D = 7.4 },
    @{ Row = 21; C = "01/01/2024"; D = 7.9 },
    @{ Row = 22; C = "01/04/2024"; D = 6.9 },
    @{ Row = 23; C = "01/04/2019"; D = 14.8 },
    @{ Row = 24; C = "01/07/2019"; D = 14.6 },
    @{ Row = 25; C = "01/10/2019"; D = 13.8 },
    @{ Row = 26; C = "01/01/2020"; D = 15.8 },
    @{ Row = 27; C = "01/04/2020"; D = $null },
    @{ Row = 28; C = "01/07/2020"; D = $null },
    @{ Row = 29; C = "01/10/2020"; D = $null },
    @{ Row = 30; C = "01/01/2021"; D = $null },
    @{ Row = 31; C = "01/04/2021"; D = $null },
    @{ Row = 32; C = "01/07/2021"; D = $null },
    @{ Row = 33; C = "01/10/2021"; D = $null },
    @{ Row = 34; C = "01/01/2022"; D = $null },
    @{ Row = 35; C = "01/04/2022"; D = 12.7 },
    @{ Row = 36; C = "01/07/2022"; D = 12 },
    @{ Row = 37; C = "01/10/2022"; D = 10.9 },
    @{ Row = 38; C = "01/01/2023"; D = 12.2 },
    @{ Row = 39; C = "01/04/2023"; D = 11.3 },
    @{ Row = 40; C = "01/07/2023"; D = 10.8 },
    @{ Row = 41; C = "01/10/2023"; D = 10.4 },
    @{ Row = 42; C = "01/01/2024"; D = 11.1 },
    @{ Row = 43; C = "01/04/2024"; D = 9.4 },
    @{ Row = 44; C = "01/04/2019"; D = 15.4 },
    @{ Row = 45; C = "01/07/2019"; D = 14.8 },
    @{ Row = 46; C = "01/10/2019"; D = 15 },
    @{ Row = 47; C = "01/01/2020"; D = 15.8 },
    @{ Row = 48; C = "01/04/2020"; D = $null },
    @{ Row = 49; C = "01/07/2020"; D = $null },
    @{ Row = 50; C = "01/10/2020"; D = $null },
    @{ Row = 51; C = "01/01/2021"; D = $null },
    @{ Row = 52; C = "01/04/2021"; D = $null },
    @{ Row = 53; C = "01/07/2021"; D = $null },
    @{ Row = 54; C = "01/10/2021"; D = $null },
    @{ Row = 55; C = "01/01/2022"; D = $null },
    @{ Row = 56; C = "01/04/2022"; D = 12.7 },
    @{ Row = 57; C = "01/07/2022"; D = 12.1 },
    @{ Row = 58; C = "01/10/2022"; D = 11.9 },
    @{ Row = 59; C = "01/01/2023"; D = 11.9 },
    @{ Row = 60; C = "01/04/2023"; D = 10.3 },
    @{ Row = 61; C = "01/07/2023"; D = 9.800000000000001 },
    @{ Row = 62; C = "01/10/2023"; D = 11.2 },
    @{ Row = 63; C = "01/01/2024"; D = 10 },
    @{ Row = 64; C = "01/04/2024"; D = 9.1 }
)

foreach ($u in $updates) {
    # Force the date-like text into column C as plain text (NumberFormat "@")
    # so Excel doesn't silently reinterpret "01/04/2019" as a date serial;
    # resetting the Style afterwards avoids leaving a stray text format applied.
    $cCell = $ws.Cells.Item($u.Row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C
    $cCell.Style = "Normal"

    $dCell = $ws.Cells.Item($u.Row, 4)
    if ($null -eq $u.D) {
        $dCell.ClearContents()
    } else {
        $dCell.Value = $u.D
    }
}
